$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (edit in place via Characters to preserve run position) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "43"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "10/20/2025"
$c9.Characters(48, 10).Text = "10/26/2025"

# --- Column width updates (cols I & J widen to match col H) ---
$refWidth = $ws.Columns.Item(8).ColumnWidth
$ws.Columns.Item(9).ColumnWidth = $refWidth
$ws.Columns.Item(10).ColumnWidth = $refWidth

# --- Cell value updates ---
# Same-style numeric/text updates
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 26
$ws.Range("K15").Value = 73.333333333333
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = -7.142857142857
$ws.Range("N15").Value = -44.680851063829
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 36.363636363636
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = 20.754716981132
$ws.Range("L16").Value = -0.775193798449
$ws.Range("M16").Value = -46.443514644351
$ws.Range("N16").Value = -84.466019417475
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -7.692307692307
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = -26.315789473684
$ws.Range("I17").Value = 358
$ws.Range("J17").Value = 407
$ws.Range("K17").Value = -12.039312039312
$ws.Range("L17").Value = -7.253886010362
$ws.Range("M17").Value = 44.354838709677
$ws.Range("N17").Value = -45.675265553869
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 250
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -18.691588785046
$ws.Range("M18").Value = -66.015625
$ws.Range("N18").Value = -93.937282229965
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 20.689655172413
$ws.Range("I19").Value = 290
$ws.Range("J19").Value = 323
$ws.Range("K19").Value = -10.216718266253
$ws.Range("L19").Value = -16.666666666666
$ws.Range("M19").Value = -13.946587537092
$ws.Range("N19").Value = -39.075630252100
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 11.111111111111
$ws.Range("I20").Value = 68
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = -37.614678899082
$ws.Range("L20").Value = -48.091603053435
$ws.Range("M20").Value = -59.281437125748
$ws.Range("N20").Value = -93.784277879341
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 17.857142857142
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = 16.304347826087
$ws.Range("I21").Value = 960
$ws.Range("J21").Value = 1048
$ws.Range("K21").Value = -8.396946564885
$ws.Range("L21").Value = -14.742451154529
$ws.Range("M21").Value = -25.349922239502
$ws.Range("N21").Value = -78.901098901098
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 60
$ws.Range("I23").Value = 59
$ws.Range("J23").Value = 68
$ws.Range("K23").Value = -13.235294117647
$ws.Range("L23").Value = -16.901408450704
$ws.Range("M23").Value = 55.263157894736
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -6.666666666666
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = -3.409090909090
$ws.Range("I24").Value = 987
$ws.Range("J24").Value = 1005
$ws.Range("K24").Value = -1.791044776119
$ws.Range("L24").Value = -3.894839337877
$ws.Range("M24").Value = -21.542130365659
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 3.448275862068
$ws.Range("I25").Value = 477
$ws.Range("J25").Value = 434
$ws.Range("K25").Value = 9.907834101382
$ws.Range("L25").Value = 23.255813953488
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 65
$ws.Range("G26").Value = 82
$ws.Range("H26").Value = -20.731707317073
$ws.Range("I26").Value = 644
$ws.Range("J26").Value = 658
$ws.Range("K26").Value = -2.127659574468
$ws.Range("L26").Value = 7.154742096505
$ws.Range("M26").Value = -39.130434782608
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 37
$ws.Range("K27").Value = 23.333333333333
$ws.Range("L27").Value = 94.736842105263
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 74
$ws.Range("J28").Value = 73
$ws.Range("K28").Value = 1.369863013698
$ws.Range("L28").Value = -3.896103896103
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = -30
$ws.Range("J30").Value = 10
$ws.Range("K30").Value = -30
$ws.Range("L33").Value = -25

# Style-transition updates (value + style copy from stable donor cells)
$ws.Range("D15").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C20").Value = 5
$ws.Range("F14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C29").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("F14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = 0
$ws.Range("K14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("C30").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("F14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = 0
$ws.Range("K14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
